$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# D-column cells must stay plain text (prices use "." as a thousands
# separator in some rows, e.g. "30.494.32", and some values have
# significant trailing zeros, e.g. "0.07560") so we force Text format
# before writing, then drop back to the default "Normal" style so no
# extra formatting is left behind on the cell.

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "30.494.32"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +0.17%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "1.914.50"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  -0.11%  "

$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9985"
$dCell.Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "244.80"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  +0.45%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9987"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  -0.14%  "

$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "0.4797"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  +2.19%  "

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "0.2884"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  +0.57%  "

$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06717"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  -1.82%  "

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "110.66"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  +0.10%  "

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "19.23"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  +4.38%  "

$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "1.913.53"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  -0.09%  "

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "0.07560"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  -2.16%  "

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "5.232"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  -1.03%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "0.6660"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  +1.38%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "303.28"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  +2.15%  "

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "30.465.70"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("E18").Value = "  -0.02%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9981"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  -0.21%  "

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "0.000007562"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  -0.93%  "

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "2.156.14"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  +0.50%  "

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "5.469"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  +4.33%  "

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9984"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  -0.15%  "

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "6.389"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  +2.85%  "

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "9.466"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  +0.97%  "

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "164.37"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  -2.66%  "

$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "20.52"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  -5.77%  "

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "2.086"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  -0.46%  "

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "0.1070"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("E30").Value = "  +2.52%  "

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "4.164"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  -0.28%  "

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "4.017"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  +0.89%  "

$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "0.04974"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  -1.22%  "

$ws.Range("E34").Value = "  -0.28%  "

$ws.Range("E35").Value = "  -1.41%  "

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9986"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("E37").Value = "  -1.26%  "

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "2.726"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  -0.51%  "

$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "2.673"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  -0.13%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "111.39"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  +1.69%  "

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "2.016"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  -2.03%  "

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "0.4407"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  +3.66%  "

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "0.8647"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  -0.63%  "

$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "5.907"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  +0.95%  "

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "0.9985"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  -0.14%  "

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "68.53"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  +1.78%  "

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "49.64"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  -3.81%  "

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "7.272"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  +1.09%  "

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "9.278"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  +0.32%  "

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "0.1232"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  +1.16%  "

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "0.2537"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  +3.88%  "
